# Auto-applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '28.082.74'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +1.11%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.817.00'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +1.70%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9977'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.35%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '310.19'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9991'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.19%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5010'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -1.94%  '

# Row 8
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +1.54%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.09963'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +27.28%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.106'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  +1.54%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '40.83'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '6.429'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +3.39%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '20.57'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  +1.99%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.9972'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -0.38%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '1.808.65'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +2.01%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '7.285'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +0.91%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001141'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  +6.06%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '92.50'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +1.41%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06645'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +1.78%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.9991'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -0.21%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.20'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +1.13%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.933'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '28.118.66'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '11.12'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +0.97%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.257'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +1.58%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '158.75'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -1.24%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '20.63'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +2.03%  '

# Row 28
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.018.19'
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +1.54%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.417'
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +2.51%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '127.29'
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +3.04%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.1066'
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.77%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.037'
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +0.28%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.582'
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +1.84%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.600'
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -0.80%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.06745'
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -4.74%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '8.949'
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +2.38%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.02338'
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +1.44%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.2143'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +0.99%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '4.956'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -0.58%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '11.31'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -1.39%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.6208'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +2.11%  '

# Row 42
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +2.59%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '13.17'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +0.18%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.5921'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +0.67%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.694'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.281'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -2.40%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '124.49'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +0.37%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.930'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +1.52%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.180'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.06786'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -0.37%  '
